$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = '@'
$c.Value = '62.713.31'
$c = $ws.Range("E2")
$c.NumberFormat = '@'
$c.Value = '  -0.64%  '
$c = $ws.Range("D3")
$c.NumberFormat = '@'
$c.Value = '2.580.08'
$c = $ws.Range("E3")
$c.NumberFormat = '@'
$c.Value = '  +0.74%  '
$c = $ws.Range("E4")
$c.NumberFormat = '@'
$c.Value = '  +0.03%  '
$c = $ws.Range("D5")
$c.NumberFormat = '@'
$c.Value = '580.45'
$c = $ws.Range("E5")
$c.NumberFormat = '@'
$c.Value = '  -0.24%  '
$c = $ws.Range("D6")
$c.NumberFormat = '@'
$c.Value = '144.16'
$c = $ws.Range("E6")
$c.NumberFormat = '@'
$c.Value = '  -2.36%  '
$c = $ws.Range("E7")
$c.NumberFormat = '@'
$c.Value = '  +0.01%  '
$c = $ws.Range("E8")
$c.NumberFormat = '@'
$c.Value = '  +0.71%  '
$c = $ws.Range("E9")
$c.NumberFormat = '@'
$c.Value = '  -0.36%  '
$c = $ws.Range("E10")
$c.NumberFormat = '@'
$c.Value = '  -0.40%  '
$c = $ws.Range("E11")
$c.NumberFormat = '@'
$c.Value = '  -0.58%  '
$c = $ws.Range("E12")
$c.NumberFormat = '@'
$c.Value = '  -0.80%  '
$c = $ws.Range("D13")
$c.NumberFormat = '@'
$c.Value = '26.92'
$c = $ws.Range("E13")
$c.NumberFormat = '@'
$c.Value = '  -2.26%  '
$c = $ws.Range("D14")
$c.NumberFormat = '@'
$c.Value = '3.043.78'
$c = $ws.Range("E14")
$c.NumberFormat = '@'
$c.Value = '  +0.75%  '
$c = $ws.Range("D15")
$c.NumberFormat = '@'
$c.Value = '62.646.51'
$c = $ws.Range("E15")
$c.NumberFormat = '@'
$c.Value = '  -0.62%  '
$c = $ws.Range("D16")
$c.NumberFormat = '@'
$c.Value = '0.0000144'
$c = $ws.Range("E16")
$c.NumberFormat = '@'
$c.Value = '  -0.20%  '
$c = $ws.Range("D17")
$c.NumberFormat = '@'
$c.Value = '2.596.78'
$c = $ws.Range("E17")
$c.NumberFormat = '@'
$c.Value = '  +1.54%  '
$c = $ws.Range("D18")
$c.NumberFormat = '@'
$c.Value = '11.20'
$c = $ws.Range("E18")
$c.NumberFormat = '@'
$c.Value = '  -1.45%  '
$c = $ws.Range("D19")
$c.NumberFormat = '@'
$c.Value = '338.80'
$c = $ws.Range("E19")
$c.NumberFormat = '@'
$c.Value = '  -0.91%  '
$c = $ws.Range("E20")
$c.NumberFormat = '@'
$c.Value = '  -0.54%  '
$c = $ws.Range("D21")
$c.NumberFormat = '@'
$c.Value = '6.67'
$c = $ws.Range("E21")
$c.NumberFormat = '@'
$c.Value = '  -2.13%  '
$c = $ws.Range("E22")
$c.NumberFormat = '@'
$c.Value = '  +0.09%  '
$c = $ws.Range("D23")
$c.NumberFormat = '@'
$c.Value = '67.39'
$c = $ws.Range("E23")
$c.NumberFormat = '@'
$c.Value = '  +2.24%  '
$c = $ws.Range("D24")
$c.NumberFormat = '@'
$c.Value = '1.55'
$c = $ws.Range("E24")
$c.NumberFormat = '@'
$c.Value = '  +5.11%  '
$c = $ws.Range("E25")
$c.NumberFormat = '@'
$c.Value = '  -2.62%  '
$c = $ws.Range("E26")
$c.NumberFormat = '@'
$c.Value = '  -3.02%  '
$c = $ws.Range("D27")
$c.NumberFormat = '@'
$c.Value = '8.00'
$c = $ws.Range("E27")
$c.NumberFormat = '@'
$c.Value = '  +1.33%  '
$c = $ws.Range("D28")
$c.NumberFormat = '@'
$c.Value = '0.997'
$c = $ws.Range("E28")
$c.NumberFormat = '@'
$c.Value = '  -0.32%  '
$c = $ws.Range("D29")
$c.NumberFormat = '@'
$c.Value = '8.24'
$c = $ws.Range("E29")
$c.NumberFormat = '@'
$c.Value = '  -2.20%  '
$c = $ws.Range("D30")
$c.NumberFormat = '@'
$c.Value = '1.94'
$c = $ws.Range("E30")
$c.NumberFormat = '@'
$c.Value = '  -0.75%  '
$c = $ws.Range("D31")
$c.NumberFormat = '@'
$c.Value = '0.0₃0806'
$c = $ws.Range("E31")
$c.NumberFormat = '@'
$c.Value = '  -1.93%  '
$c = $ws.Range("D32")
$c.NumberFormat = '@'
$c.Value = '456.31'
$c = $ws.Range("E32")
$c.NumberFormat = '@'
$c.Value = '  +5.92%  '
$c = $ws.Range("D33")
$c.NumberFormat = '@'
$c.Value = '176.37'
$c = $ws.Range("E33")
$c.NumberFormat = '@'
$c.Value = '  -0.04%  '
$c = $ws.Range("E34")
$c.NumberFormat = '@'
$c.Value = '  +1.29%  '
$c = $ws.Range("E35")
$c.NumberFormat = '@'
$c.Value = '  +0.18%  '
$c = $ws.Range("E36")
$c.NumberFormat = '@'
$c.Value = '  -1.55%  '
$c = $ws.Range("D37")
$c.NumberFormat = '@'
$c.Value = '18.93'
$c = $ws.Range("E37")
$c.NumberFormat = '@'
$c.Value = '  -1.45%  '
$c = $ws.Range("E38")
$c.NumberFormat = '@'
$c.Value = '  -0.62%  '
$c = $ws.Range("D40")
$c.NumberFormat = '@'
$c.Value = '1.68'
$c = $ws.Range("E40")
$c.NumberFormat = '@'
$c.Value = '  -3.28%  '
$c = $ws.Range("D41")
$c.NumberFormat = '@'
$c.Value = '159.48'
$c = $ws.Range("E41")
$c.NumberFormat = '@'
$c.Value = '  +5.01%  '
$c = $ws.Range("E42")
$c.NumberFormat = '@'
$c.Value = '  -2.65%  '
$c = $ws.Range("E43")
$c.NumberFormat = '@'
$c.Value = '  +4.08%  '
$c = $ws.Range("D44")
$c.NumberFormat = '@'
$c.Value = '21.00'
$c = $ws.Range("E44")
$c.NumberFormat = '@'
$c.Value = '  -0.13%  '
$c = $ws.Range("E45")
$c.NumberFormat = '@'
$c.Value = '  -2.96%  '
$c = $ws.Range("E46")
$c.NumberFormat = '@'
$c.Value = '  -1.05%  '
$c = $ws.Range("D47")
$c.NumberFormat = '@'
$c.Value = '0.0234'
$c = $ws.Range("E47")
$c.NumberFormat = '@'
$c.Value = '  -3.06%  '
$c = $ws.Range("E48")
$c.NumberFormat = '@'
$c.Value = '  -1.47%  '
$c = $ws.Range("E49")
$c.NumberFormat = '@'
$c.Value = '  +0.41%  '
$c = $ws.Range("E50")
$c.NumberFormat = '@'
$c.Value = '  -1.97%  '
$c = $ws.Range("E51")
$c.NumberFormat = '@'
$c.Value = '  +3.99%  '
